# Normalize the "Recorded By" (column G) values so that the "System" token
# (case-insensitive) is no longer the first entry in the comma-separated list.
# Rule observed from the source diff:
#   - If the list has no "System" entry at all, swap the first and last entries.
#   - If the list's first entry is "System" (any case), swap the first and last
#     entries (this moves "System" to the end while preserving its own casing
#     and the casing of whichever entry was last).
#   - Otherwise (a "System" entry exists but is not first) leave the value as is.
#   - Single-entry values are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $original = $cell.Text

    if ([string]::IsNullOrEmpty($original)) {
        continue
    }

    $parts = $original -split ", "

    if ($parts.Length -lt 2) {
        continue
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $hasSystem = $true
        }
    }

    $firstIsSystem = ($parts[0].ToLower() -eq "system")

    if ((-not $hasSystem) -or $firstIsSystem) {
        $lastIndex = $parts.Length - 1
        $temp = $parts[0]
        $parts[0] = $parts[$lastIndex]
        $parts[$lastIndex] = $temp
        $newValue = [string]::Join(", ", $parts)
        $ws.Cells.Item($r, 7).Value = $newValue
    }
}
